$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column C (rows 2-13), produced by re-running the
# computation loop (yields slightly different numbers than before).
$newValues = @{
    2  = -5.143960429218957
    3  = -1.225135039214685
    4  = -0.07030676086051522
    5  = -0.4261225642150345
    6  = 0.01177441783374552
    7  = 0.09949832726891362
    8  = 0.1241200437525787
    9  = 0.0292001024991827
    10 = 0.02926903675507106
    11 = 0.002762913964716858
    12 = 0.04060450717413161
    13 = -0.001865287044813922
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}
